$d = $word.ActiveDocument

# 1. Merge the two runs "多云，今天是六一" + bookmark + "儿童节，又是开心的一天呢"
#    into a single run, dropping the _GoBack bookmark that sat between them.
$find = $d.Content.Find
$find.Execute("多云，今天是六一儿童节，又是开心的一天呢", $true, $false, $false, $false, $false, $true, 1, $false, `
    "多云，今天是六一儿童节，又是开心的一天呢", 2) | Out-Null

# 2. Split "2022年6月3日星期五" into "2022年6月3日星期" + "五", with a
#    _GoBack bookmark inserted right before the final character.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "2022年6月3日星期五`r") {
        $pEnd = $p.Range.End
        $bmRange = $d.Range($pEnd - 2, $pEnd - 2)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}

# 3. Replace the text of the last paragraph (which currently reads
#    "中雨，今天是农历五月初五，中国传统端午节。") with the new exam-day text,
#    before any new paragraphs are inserted, so the Find match is unambiguous.
$find2 = $d.Content.Find
$find2.Execute("中雨，今天是农历五月初五，中国传统端午节。", $true, $false, $false, $false, $false, $true, 1, $false, `
    "晴，今天是高考的一天，上午考语文，下午考数学。", 2) | Out-Null

# 4. Insert two new paragraphs after the "2022年6月3日星期...五" paragraph:
#    "中雨，今天是农历五月初五，中国传统端午节。" and "2022年6月7日星期二"
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith("2022年6月3日星期")) {
        $insertAfter = $p.Range
        $insertAfter.Collapse(0)
        $insertAfter.InsertParagraphAfter()
        $newPara1 = $d.Paragraphs.Item($i + 1)
        $newPara1.Range.InsertBefore("中雨，今天是农历五月初五，中国传统端午节。")

        $newPara1 = $d.Paragraphs.Item($i + 1)
        $afterNew1 = $newPara1.Range
        $afterNew1.Collapse(0)
        $afterNew1.InsertParagraphAfter()
        $newPara2 = $d.Paragraphs.Item($i + 2)
        $newPara2.Range.InsertBefore("2022年6月7日星期二")
        break
    }
}
